# Apply balance-sheet value corrections to CALIFORNIA_STATE_UNIVERSITY sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated figures
$ws.Range("B3").Value = 5793096
$ws.Range("B6").Value = 2537879
$ws.Range("B9").Value = 9801435
$ws.Range("B10").Value = 22820815

# These two figures are no longer available -- clear them out entirely.
$ws.Range("B15").ClearContents()
$ws.Range("B17").ClearContents()
